# LOM3254.xlsx update
# - Rewrites the "Objetivos" answer (B10/C10) to the short text
#   "5982760 - Carlos Alberto Baldan".
# - Rebuilds rows 13-24 (previously 13-25) with a new label/value layout
#   (several rows shifted/re-paired, row 25 removed entirely).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($Col, $Row, $Text) {
    $ref = "$Col$Row"
    # Pull correctly-styled formatting for this column from an untouched
    # template row (row 3 carries the plain A/B/C styles 1/2/3) so freshly
    # created cells match the workbook's existing look.
    $tmpl = "$Col" + "3"
    $ws.Range($tmpl).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $Text
}

# --- Objetivos answer shrinks to a single short line -----------------------
$ws.Range("B10").Value = "5982760 - Carlos Alberto Baldan"
$ws.Range("C10").Value = "5982760 - Carlos Alberto Baldan"

# --- Wipe the old rows 13-25 and rebuild 13-24 from scratch -----------------
$ws.Range("A13:C25").EntireRow.Delete()

Set-Cell "A" 13 "Programa resumido:"
Set-Cell "B" 13 "Semestral"
Set-Cell "C" 13 "Semestral"
$ws.Rows.Item(13).RowHeight = 60

Set-Cell "A" 14 "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

Set-Cell "A" 15 "Programa:"
Set-Cell "B" 15 "01/01/2015"
Set-Cell "C" 15 "01/01/2015"
$ws.Rows.Item(15).RowHeight = 120

Set-Cell "A" 16 "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

Set-Cell "A" 17 "Avaliação:"

Set-Cell "A" 18 "Método:"
Set-Cell "B" 18 "5982760 - Carlos Alberto Baldan"
Set-Cell "C" 18 "5982760 - Carlos Alberto Baldan"
$ws.Rows.Item(18).RowHeight = 60

Set-Cell "A" 19 "Critério:"
Set-Cell "B" 19 "Aulas práticas utilizando placas contendo circuitos elétricos para medições com equipamentos eletrônicos. Máximo 2 alunos por bancada e 20 alunos por turma  para 10 bancadas."
Set-Cell "C" 19 "Aulas práticas utilizando placas contendo circuitos elétricos para medições com equipamentos eletrônicos. Máximo 2 alunos por bancada e 20 alunos por turma  para 10 bancadas."
$ws.Rows.Item(19).RowHeight = 60

Set-Cell "A" 20 "Norma de recuperação:"
Set-Cell "B" 20 "Duas provas (P1 e P2) e relatórios sobre tópicos da disciplina.`nNF= MR*0,2 +0,8*(P1 + P2)/2...........MR- média de notas dos relatórios"
Set-Cell "C" 20 "Duas provas (P1 e P2) e relatórios sobre tópicos da disciplina.`nNF= MR*0,2 +0,8*(P1 + P2)/2...........MR- média de notas dos relatórios"
$ws.Rows.Item(20).RowHeight = 60

Set-Cell "A" 21 "Bibliografia:"
Set-Cell "B" 21 "(NF + RC)/2"
Set-Cell "C" 21 "(NF + RC)/2"
$ws.Rows.Item(21).RowHeight = 120

Set-Cell "A" 22 "Requisitos:"

Set-Cell "B" 23 "LOB1053 -  Física III  (Requisito)`n"
Set-Cell "C" 23 "LOB1053 -  Física III  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

Set-Cell "B" 24 "LOM3202 -  Circuitos Elétricos  (Indicação de Conjunto)`n"
Set-Cell "C" 24 "LOM3202 -  Circuitos Elétricos  (Indicação de Conjunto)`n"
$ws.Rows.Item(24).RowHeight = 30
